$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (was 45406 -> 45436, i.e. one month later)
$ws.Range("A1").Value = 45436

# Update prices in the first price table
$ws.Range("D14").Value = 271.126
$ws.Range("D15").Value = 415.87

# Update prices in the second price table
$ws.Range("D38").Value = 499.042
$ws.Range("D39").Value = 535.769
